# Auto-generated Excel COM-interop PowerShell edit script
# Implements: insert guest test-case row (Car_Details) for "Car_Details_and_Reservation.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at row 3 -- shifts existing rows 3-7 down to 4-8,
#    Excel copies row-2 formatting/styles into the new row automatically,
#    and the merged cell (old A4:K4) shifts down to A5:K5 automatically.
$ws.Rows("3:3").Insert()
$ws.Rows(3).RowHeight = 131.25

# 2) The two "car details" mockup pictures are anchored by (row, rowOff) and
#    do not auto-shift with the row insert, so nudge each down by exactly
#    one rows height so they stay pinned next to the same text rows.
$shapes = $ws.Shapes
$shp1 = $shapes.Item(1)
$shp2 = $shapes.Item(2)
$shp1.Top = $shp1.Top + $ws.Rows(3).RowHeight
$shp2.Top = $shp2.Top + $ws.Rows(3).RowHeight

# 3) Row 4 (previously row 3, "Car_Details_02" UI case): renumber its test
#    case ID to "Car_Details_03" now that it is third in the list.
$ws.Range("A4").Value = "Car_Details_03"

# 4) New row 3: guest "See more" test case.
$ws.Range("A3").Value = "Car_Details_02"
$ws.Range("B3").Value = "Car_SRS_07"
$ws.Range("C3").Value = "Functional"
$ws.Range("E3").Value = "1) open URL ""http://CarPurchasing""`n2)don't login"
$ws.Range("G3").Value = "1)From home page click on ""see more"" button at any car"
$ws.Range("H3").Value = "Guest should be redirected to a registration  page"
$ws.Range("J3").Value = "Fatma"
$ws.Range("K3").Value = "passed"
$ws.Range("D3").Value = "Validate ""See more"" button functionality for a guest"

# 5) Row 2 (Car_Details_01): retitle to clarify it is the "user" flow.
$ws.Range("D2").Value = "Validate ""See more"" button functionality for a user"

# 6) Record the reviewer ("jannat") under "Reviewed by" (M) for rows 2, 4,
#    6-8 -- previously it had been miskeyed under "Bug ID" (L) for the
#    reservation rows (now 6-8), so clear those first.
$ws.Range("M2").Value = "jannat"
$ws.Range("M4").Value = "jannat"
$ws.Range("L6").ClearContents()
$ws.Range("M6").Value = "jannat"
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = "jannat"
$ws.Range("L8").ClearContents()
$ws.Range("M8").Value = "jannat"

# 7) Restore selection/active cell as in the saved workbook.
$ws.Range("M3").Select()

